$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column E (duplicate_image_filename) with "NA" for data rows 2 through 21.
$ws.Range("E2:E21").Value = "NA"
